# New weekly price observation for "Betarraga" (Vega Modelo de Temuco) is
# inserted as row 372, pushing the existing rows 372..429 down to 373..430.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(372).Insert()

$ws.Range("A372").Value = 10
$ws.Range("B372").Value = "Vega Modelo de Temuco"
$ws.Range("C372").Value = "La Araucanía"
$ws.Range("D372").Value = 44776
$ws.Range("E372").Value = 9
$ws.Range("F372").Value = 100114014
$ws.Range("G372").Value = "Betarraga"
$ws.Range("H372").Value = "Sin especificar"
$ws.Range("I372").Value = "Primera"
$ws.Range("J372").Value = 100
$ws.Range("K372").Value = 9000
$ws.Range("L372").Value = 9000
$ws.Range("M372").Value = 9000
$ws.Range("N372").Value = '$/docena de paquetes'
$ws.Range("O372").Value = "Región del Maule"
$ws.Range("P372").Value = 750
$ws.Range("Q372").Value = 12
$ws.Range("R372").Value = "Hortaliza"
